$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474758173/gFRHQATQ91-camera-video-segment-driver-1749474760673.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK3AXEIZ%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T170045Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIGsK87waLbg%2F38%2Bq22%2F5NWGDSjjYjcPgCcyR6GZ3WBZ0AiAmJiH5MwtqvR%2FmFF6fSZzRQTsTevRt96%2BjPQYnqHpjVCrmAwjB%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMsLyKjggqJGtzubB1KroDGMzjqXHYybvGU6uHL6%2Fx7v%2Fy1%2BSvIBOjGMZrT%2BNzT3KdkY2zheG5LmRFGLyR88z7Zfk%2FN%2FnMTYEv8NWUcdNlTVUcqIw9iooDIwkG2Mvlt1AwJInE7mbbjg7Ly86y%2Bq43vF%2BLPyWwCF8GXXSFgbkhq4lb2cwHdLy2sffgRJsKygRQUVfIQH5PknuUvv4k1T5nrPG0JSRm6Ogt%2Bs7HZ2cTWNfx1%2BdFwmw2JTOIWc0eNTItWT%2Bc6aZl2KcG1roGatXzzCwGGyxFBDmKLPJtOyePy%2FZ5rz1oYGJDBRqsDyd3vi3xzCMQOqZkkOguOFbQgwjRXvKrDj7YL5H9HfGSEFViwL3FbByXYjitheQwgzMfB1rhmiYUzaJYYsLscYWrhsaD1YMLtTtw9tjRZrGJtGNe64cP8o21TU3JfTZ%2FI%2B6rZq6feRbXl758H1N8HDqAHuy2I4ab27R%2BsUvuVobnRJa%2FrRqgl9c71R59j%2BAavfzf9c%2FpjhNZthVkNtm3TIxOvJtd7%2B79%2FCT6qmAX3nDEQkExwC%2BQNi%2BtB9lGMMbg%2BNibz6f8zDnvSSvJ0aTifn%2BB9pZzrPYCmSJRc9WNFzDoqaHCBjqmAZ8GS4GtlF8wttkw6h4hQp9eDUHmbaMN%2Bkc8VZOW4t3HXvkX9jHZ3yPMfI3d5gCvhuHPolTSRBu6ziASJOC7aPMeifU8Kkt7a91esXEJtf%2FjzQVqjzABBdYqOGzTTwHSe%2BrvI0eYfhvRvGRA36MHXjhoQJfddkx%2BvKyJuU4EgVBNHIowE49GhlUhg99AJIGGgcSbc1wHQISNkaRVlnGg%2FXkaepm0fMQ%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2001%3A00%3A45%20GMT&X-Amz-Signature=fe84bd91c3724aedf99d520593c718251c466857766b59eb9a751c0fc83278b7'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474758173/7Z3ova6Wtr-camera-video-segment-1749474760673.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK3AXEIZ%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T170045Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIGsK87waLbg%2F38%2Bq22%2F5NWGDSjjYjcPgCcyR6GZ3WBZ0AiAmJiH5MwtqvR%2FmFF6fSZzRQTsTevRt96%2BjPQYnqHpjVCrmAwjB%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMsLyKjggqJGtzubB1KroDGMzjqXHYybvGU6uHL6%2Fx7v%2Fy1%2BSvIBOjGMZrT%2BNzT3KdkY2zheG5LmRFGLyR88z7Zfk%2FN%2FnMTYEv8NWUcdNlTVUcqIw9iooDIwkG2Mvlt1AwJInE7mbbjg7Ly86y%2Bq43vF%2BLPyWwCF8GXXSFgbkhq4lb2cwHdLy2sffgRJsKygRQUVfIQH5PknuUvv4k1T5nrPG0JSRm6Ogt%2Bs7HZ2cTWNfx1%2BdFwmw2JTOIWc0eNTItWT%2Bc6aZl2KcG1roGatXzzCwGGyxFBDmKLPJtOyePy%2FZ5rz1oYGJDBRqsDyd3vi3xzCMQOqZkkOguOFbQgwjRXvKrDj7YL5H9HfGSEFViwL3FbByXYjitheQwgzMfB1rhmiYUzaJYYsLscYWrhsaD1YMLtTtw9tjRZrGJtGNe64cP8o21TU3JfTZ%2FI%2B6rZq6feRbXl758H1N8HDqAHuy2I4ab27R%2BsUvuVobnRJa%2FrRqgl9c71R59j%2BAavfzf9c%2FpjhNZthVkNtm3TIxOvJtd7%2B79%2FCT6qmAX3nDEQkExwC%2BQNi%2BtB9lGMMbg%2BNibz6f8zDnvSSvJ0aTifn%2BB9pZzrPYCmSJRc9WNFzDoqaHCBjqmAZ8GS4GtlF8wttkw6h4hQp9eDUHmbaMN%2Bkc8VZOW4t3HXvkX9jHZ3yPMfI3d5gCvhuHPolTSRBu6ziASJOC7aPMeifU8Kkt7a91esXEJtf%2FjzQVqjzABBdYqOGzTTwHSe%2BrvI0eYfhvRvGRA36MHXjhoQJfddkx%2BvKyJuU4EgVBNHIowE49GhlUhg99AJIGGgcSbc1wHQISNkaRVlnGg%2FXkaepm0fMQ%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2001%3A00%3A45%20GMT&X-Amz-Signature=9cee993b4292ff3a82b4d372830794a391c76df5cfde909da2886115bc7ce3b4'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474688613/JLWgs3viB6-camera-video-segment-driver-1749474691113.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK3AXEIZ%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T170045Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIGsK87waLbg%2F38%2Bq22%2F5NWGDSjjYjcPgCcyR6GZ3WBZ0AiAmJiH5MwtqvR%2FmFF6fSZzRQTsTevRt96%2BjPQYnqHpjVCrmAwjB%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMsLyKjggqJGtzubB1KroDGMzjqXHYybvGU6uHL6%2Fx7v%2Fy1%2BSvIBOjGMZrT%2BNzT3KdkY2zheG5LmRFGLyR88z7Zfk%2FN%2FnMTYEv8NWUcdNlTVUcqIw9iooDIwkG2Mvlt1AwJInE7mbbjg7Ly86y%2Bq43vF%2BLPyWwCF8GXXSFgbkhq4lb2cwHdLy2sffgRJsKygRQUVfIQH5PknuUvv4k1T5nrPG0JSRm6Ogt%2Bs7HZ2cTWNfx1%2BdFwmw2JTOIWc0eNTItWT%2Bc6aZl2KcG1roGatXzzCwGGyxFBDmKLPJtOyePy%2FZ5rz1oYGJDBRqsDyd3vi3xzCMQOqZkkOguOFbQgwjRXvKrDj7YL5H9HfGSEFViwL3FbByXYjitheQwgzMfB1rhmiYUzaJYYsLscYWrhsaD1YMLtTtw9tjRZrGJtGNe64cP8o21TU3JfTZ%2FI%2B6rZq6feRbXl758H1N8HDqAHuy2I4ab27R%2BsUvuVobnRJa%2FrRqgl9c71R59j%2BAavfzf9c%2FpjhNZthVkNtm3TIxOvJtd7%2B79%2FCT6qmAX3nDEQkExwC%2BQNi%2BtB9lGMMbg%2BNibz6f8zDnvSSvJ0aTifn%2BB9pZzrPYCmSJRc9WNFzDoqaHCBjqmAZ8GS4GtlF8wttkw6h4hQp9eDUHmbaMN%2Bkc8VZOW4t3HXvkX9jHZ3yPMfI3d5gCvhuHPolTSRBu6ziASJOC7aPMeifU8Kkt7a91esXEJtf%2FjzQVqjzABBdYqOGzTTwHSe%2BrvI0eYfhvRvGRA36MHXjhoQJfddkx%2BvKyJuU4EgVBNHIowE49GhlUhg99AJIGGgcSbc1wHQISNkaRVlnGg%2FXkaepm0fMQ%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2001%3A00%3A45%20GMT&X-Amz-Signature=b65cb1a7ec1fa2fb5f5c1bcc25442eae1fcdf59e4d4ada5baffb601d0ccd77bd'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474688613/LOHXs2uAgt-camera-video-segment-1749474691113.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK3AXEIZ%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T170045Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIGsK87waLbg%2F38%2Bq22%2F5NWGDSjjYjcPgCcyR6GZ3WBZ0AiAmJiH5MwtqvR%2FmFF6fSZzRQTsTevRt96%2BjPQYnqHpjVCrmAwjB%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMsLyKjggqJGtzubB1KroDGMzjqXHYybvGU6uHL6%2Fx7v%2Fy1%2BSvIBOjGMZrT%2BNzT3KdkY2zheG5LmRFGLyR88z7Zfk%2FN%2FnMTYEv8NWUcdNlTVUcqIw9iooDIwkG2Mvlt1AwJInE7mbbjg7Ly86y%2Bq43vF%2BLPyWwCF8GXXSFgbkhq4lb2cwHdLy2sffgRJsKygRQUVfIQH5PknuUvv4k1T5nrPG0JSRm6Ogt%2Bs7HZ2cTWNfx1%2BdFwmw2JTOIWc0eNTItWT%2Bc6aZl2KcG1roGatXzzCwGGyxFBDmKLPJtOyePy%2FZ5rz1oYGJDBRqsDyd3vi3xzCMQOqZkkOguOFbQgwjRXvKrDj7YL5H9HfGSEFViwL3FbByXYjitheQwgzMfB1rhmiYUzaJYYsLscYWrhsaD1YMLtTtw9tjRZrGJtGNe64cP8o21TU3JfTZ%2FI%2B6rZq6feRbXl758H1N8HDqAHuy2I4ab27R%2BsUvuVobnRJa%2FrRqgl9c71R59j%2BAavfzf9c%2FpjhNZthVkNtm3TIxOvJtd7%2B79%2FCT6qmAX3nDEQkExwC%2BQNi%2BtB9lGMMbg%2BNibz6f8zDnvSSvJ0aTifn%2BB9pZzrPYCmSJRc9WNFzDoqaHCBjqmAZ8GS4GtlF8wttkw6h4hQp9eDUHmbaMN%2Bkc8VZOW4t3HXvkX9jHZ3yPMfI3d5gCvhuHPolTSRBu6ziASJOC7aPMeifU8Kkt7a91esXEJtf%2FjzQVqjzABBdYqOGzTTwHSe%2BrvI0eYfhvRvGRA36MHXjhoQJfddkx%2BvKyJuU4EgVBNHIowE49GhlUhg99AJIGGgcSbc1wHQISNkaRVlnGg%2FXkaepm0fMQ%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2001%3A00%3A45%20GMT&X-Amz-Signature=d68b493e960de535cbb5335e3f9d2928a164ddacc73aa629bfc41b0b36df69cf'
